$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to stay text ("@") while we assign numeric-looking strings,
# then restore the Normal style so no stray style index is left on the cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.205.98"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "3.512.76"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "595.43"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "172.38"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  +6.63%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "4.122.96"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "28.65"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "67.169.86"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "3.515.17"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "14.21"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "396.26"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "73.35"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  -5.50%  "
$ws.Range("D26").Value = "10.27"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "6.27"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "24.07"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").Value = "7.39"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "1.64"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("D35").Value = "164.07"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").Value = "0.892"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").Value = "6.89"
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("D39").Value = "4.72"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").Value = "0.0745"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "26.36"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "2.822.52"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "2.61"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").Value = "42.92"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "0.0308"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").Value = "341.79"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "33.95"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  -0.29%  "

$dRange.Style = "Normal"
